$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = "0959e04b3aa73c7868f6602bf922730d"
$ws.Range("B11").Value = "b2c8390815ce162dfbc195a6e3539e5e"
$ws.Range("B15").Value = "23c4bc40f0a8eb34227b73fcade17c44"
$ws.Range("B17").Value = "9a0cf9cde071af21b9a8b1250544dbe1"
$ws.Range("B24").Value = "97658fa5e114113e1a449acdf95ddf5c"
$ws.Range("B29").Value = "880169c4f20521e3d4822a10de244c87"
$ws.Range("B34").Value = "3f54be0653ad2244272cdb4c92b66659"
$ws.Range("B121").Value = "1b616b2e73a9c56fefbc2e46caa895da"
$ws.Range("B133").Value = "219aefafdaead1e58e3487a55809ca80"
$ws.Range("B136").Value = "4d8d0cbab3ffe559b044913127f931cb"
$ws.Range("B159").Value = "258d1be4e5ce772f9c17817b83122106"
$ws.Range("B162").Value = "5ade9c4d2c6ee935e6b926f7fb9a0ce9"
$ws.Range("B169").Value = "934471d5234116c2105632f918393f08"
$ws.Range("B175").Value = "4d0c7a05dad8d06ddc754c5606b18e82"
$ws.Range("B180").Value = "3b78fbf76c5f265df55a25de18c3e2f9"
$ws.Range("B183").Value = "3e3a66cbe6076aaf0f431ff00351763d"
$ws.Range("B191").Value = "5a7741054071e5dfde5bf9e71a36d178"
$ws.Range("B198").Value = "64b254efb3909fc569555fa116472ee4"
$ws.Range("B200").Value = "77fc6691c02ede0e98ed5720035b5c68"
$ws.Range("B213").Value = "3a717925e8b0b0dcef43f46beb6facf7"
$ws.Range("B227").Value = "5df12c5655f7fb6f31c94af54215b5aa"
$ws.Range("B228").Value = "9fdefb1cd13a71ebba21891c6d2c9ee0"
$ws.Range("B232").Value = "ee3640aa2c9fca8dbcd22cc7e942fc4a"
$ws.Range("B339").Value = "3c91afa877227368cb569ee456c97b0e"
$ws.Range("B464").Value = "cafa73b84464e6ce32c8cccad7acbb7e"
$ws.Range("B465").Value = "227de680d72f57468721c27f3cc54e37"
$ws.Range("B483").Value = "8e377676ef963f85fc6cdc072adee325"
$ws.Range("B485").Value = "b8e03041b79435988ca255308392a09e"
$ws.Range("B506").Value = "74d987e2cda486e5de1a59d10854a514"
$ws.Range("B507").Value = "444c85f4b5479d65e5f444f1d33ebf48"
$ws.Range("B508").Value = "4d537e1fa995288b61de8192a7501164"
$ws.Range("B513").Value = "ad8624bb8862b0276bdeb95a68584b86"
$ws.Range("B521").Value = "b53cb95e7b1beed1711de2295117f6fb"
$ws.Range("B524").Value = "7093e1fa3dcbb0cbb3abfe84b8119398"
$ws.Range("B532").Value = "a8f9181ed491ed1e0639f790b03e4d96"
$ws.Range("B555").Value = "2913280eaeaab28ba119c5ccfd4cc4b2"
$ws.Range("B580").Value = "2e502c7addb80191a57546bebb4ca098"
$ws.Range("B624").Value = "19ad8120ef4e7fd8c61b97404cc3a38f"
$ws.Range("B635").Value = "64dc500dba2d19c1084f441cb01c798a"
$ws.Range("B637").Value = "f6a8676f79701259379a58f88f2cf0e1"
$ws.Range("B657").Value = "ea0bb9282d0b2a34cffce36bf8ed8796"
$ws.Range("B663").Value = "39ad392d778518bcc663c52f94db70b2"
$ws.Range("B673").Value = "cbb5f3ebf4381d6e4b27c30867ccb7f7"
$ws.Range("B674").Value = "ebca48fdbfb7ccaf67e04147f6865b4e"
$ws.Range("B688").Value = "15158a0991e3dad4fd94dfa5f9c8f3aa"
$ws.Range("B693").Value = "ebe45a973afff04c51d23b1b99035c84"
$ws.Range("B708").Value = "1f4e61800299458a2b76285fe27abd7a"
$ws.Range("B711").Value = "04461bccc6ab0a10df8f8af8fdc52745"
$ws.Range("B712").Value = "9866185052e14f49b301a47e90057f55"
$ws.Range("B723").Value = "b45340bd18cd2b4943af8829769651fb"
$ws.Range("B737").Value = "49281e820c63918dbaceddd9728ab270"
$ws.Range("B741").Value = "fd03ec2e714e596c6312367eb6d1c042"
$ws.Range("B750").Value = "4c1553eee3fd1eb9927e78dac8b8963e"
$ws.Range("B827").Value = "c535bd182261cc93be3c4531f608bc46"
$ws.Range("B838").Value = "71f8b444f7700ac0320c268e6589b6c9"
$ws.Range("B843").Value = "3d731832fb79f3cbf265acdce71ca60f"
$ws.Range("B862").Value = "8c360e20f2851665840633e15dbd912e"
